$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins/Losses/Ties headers in AD1:AF1, matching the style (bold,
# bordered, centered) already used by the other header cells.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-51: every row gets the same team record (Wins/Losses/Ties) values.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 66
    $ws.Cells.Item($r, 31).Value = 96
    $ws.Cells.Item($r, 32).Value = 0
}
